$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Motivation and Summary") - Content Placeholder 2 (shape 3)
# Update first bullet: drop the bullet formatting (No Bullet) and change
# "unemployment" -> "poverty"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange
$para2_1 = $tr2.Paragraphs(1,1)
$para2_1.Text = "In Chicago, is there a correlation between housing price, income, or poverty and the walkability of an area?"
$para2_1.ParagraphFormat.Bullet.Visible = 0

# ---------------------------------------------------------------------------
# Slide 6 ("Data Analysis") - Content Placeholder 2 (shape 2)
# Replace the first paragraph's text with the new sentence, split across
# three runs so "walkscore" can be marked as a flagged/unusual word.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$para6_1 = $tr6.Paragraphs(1,1)
$para6_1.Text = "The further outside the city you get the lower the "

$tr6b = $sh6.TextFrame.TextRange
$run6_1 = $tr6b.Paragraphs(1,1)
$run6_1.InsertAfter("walkscore")

$tr6c = $sh6.TextFrame.TextRange
$para6_1c = $tr6c.Paragraphs(1,1)
$para6_1c.InsertAfter(" and the higher the income")

# ---------------------------------------------------------------------------
# Slide 12 ("Post Mortem") - Content Placeholder 2 (shape 2)
# Insert a new first bullet + blank line before "Discuss any difficulties..."
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange
$para12_1 = $tr12.Paragraphs(1,1)
$para12_1.InsertBefore("There is a lot of data- especially in the census`r`r")
